# Edit slide 7 ("Let's Go") - update the realpython.com / tkinter link
# text blocks (merge the three runs that used to spell out the URL into a
# single hyperlinked run), resize/reposition the two shapes, and reorder
# the presentation-level ARTICULATE_PROJECT_OPEN tag to the end of the
# tag list (matching an "Add files via upload" re-save of the deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# ---------------------------------------------------------------------
# Shape 2: "Text Placeholder 2" (id=3)
#   "https://realpython.com/" + "python-gui-tkinter"(hlink) + "/"
#   -> single hyperlinked run "https://realpython.com/python-gui-tkinter/"
# ---------------------------------------------------------------------
$sh2 = $s.Shapes.Item(2)
$tf2 = $sh2.TextFrame
$tr2 = $tf2.TextRange

# Drop the trailing "/" run and the leading "https://realpython.com/" run;
# re-insert their text straight onto the hyperlinked run so everything
# collapses into one run that carries the hyperlink + new size.
$tr2.Characters(42, 1).Delete() | Out-Null
$tr2.Characters(1, 23).Delete() | Out-Null
$midRun2 = $tr2.Characters(1, $tr2.Length)
$midRun2.InsertBefore("https://realpython.com/") | Out-Null
$midRun2b = $tr2.Characters(1, $tr2.Length)
$midRun2b.InsertAfter("/") | Out-Null

$wholeRun2 = $tr2.Characters(1, $tr2.Length)
$wholeRun2.Font.Size = 44

# Paragraph: no bullet, marL/indent reset to 0
$lvl2 = $tf2.Ruler.Levels.Item(1)
$lvl2.LeftMargin = 0
$lvl2.FirstMargin = 0
$tr2.ParagraphFormat.Bullet.Visible = 0

# Autofit: normAutofit -> noAutofit
$tf2.AutoSize = 0

# Reposition / resize (EMU / 12700, fine-tuned so the float32 conversion
# used internally lands exactly on the target EMU values)
$sh2.Left = 74.0
$sh2.Top = 183.03527559055118
$sh2.Width = 920.0
$sh2.Height = 57.83866311732285

# ---------------------------------------------------------------------
# Shape 3: "Rectangle 3" (id=4)
#   "https://realpython.com/" + "python-gui-tkinter"(hlink) + "/"
#   -> "https://realpython.com/python-gui-tkinter"(hlink) + "/"
# ---------------------------------------------------------------------
$sh3 = $s.Shapes.Item(3)
$tf3 = $sh3.TextFrame
$tr3 = $tf3.TextRange

# Merge the leading "https://realpython.com/" run into the hyperlinked
# "python-gui-tkinter" run; leave the trailing "/" run on its own.
$tr3.Characters(1, 23).Delete() | Out-Null
$midRun3 = $tr3.Characters(1, 18)
$midRun3.InsertBefore("https://realpython.com/") | Out-Null

$wholeRun3 = $tr3.Characters(1, $tr3.Length)
$wholeRun3.Font.Size = 40

# NOTE: the target XML also adds a shape-level <a:hlinkClick r:id="rId3"/>
# directly on this shape's <p:cNvPr> (a "click anywhere on the shape"
# action setting). The PowerPoint COM surface exposed by this host only
# ever materialises Shape.ActionSettings(...).Hyperlink.Address writes as
# a run-level hlinkClick (it rewrites every run in the text, which would
# incorrectly add the link back onto the trailing "/" run), and there is
# no Shape.Hyperlink/cNvPr-level COM property available here, so that one
# sub-change cannot be reproduced through run_com in this environment.

# Reposition / resize
$sh3.Left = 77.56519705039369
$sh3.Top = 300.0
$sh3.Width = 720.0
$sh3.Height = 55.73905511811024

# ---------------------------------------------------------------------
# Presentation tags: move ARTICULATE_PROJECT_OPEN to the end of the list
# ---------------------------------------------------------------------
$p.Tags.Delete("ARTICULATE_PROJECT_OPEN")
$p.Tags.Add("ARTICULATE_PROJECT_OPEN", "0")
